# Storing workouts: rebuild Sheet1 as a clean workout table
# (Trio / Reps / Pace / Distance) with two sample rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old leading "Workout:" column -- B:E shift left to A:D,
# carrying their formatting along (incl. the numeric style that sits on
# the "5000 3000" pace cell).
$ws.Columns.Item(1).Delete()

# New third data row's values first, then the rest of row 2/3, then the
# header row -- matches how the new entries were actually typed in.
$ws.Range("C3").Value = "5000 5000"
$ws.Range("B3").Value = "2, 2"
$ws.Range("B2").Value = "1, 1"
$ws.Range("A2").Value = "1, 2, 3"
$ws.Range("A3").Value = "4, 5, 6"
$ws.Range("A1").Value = "Trio"

$ws.Range("B1").Value = "Reps"
$ws.Range("C1").Value = "Pace"
$ws.Range("D1").Value = "Distance"
$ws.Range("C2").Value = "5000 3000"
$ws.Range("D2").Value = 50
$ws.Range("D3").Value = 70

# Leave the selection on the newly added row
$ws.Range("A3").Select()
